$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.446.17"
$ws.Range("E2").Value = "  -0.42%  "
$ws.Range("D3").Value = "3.321.63"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.42"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("E7").Value = "  +8.17%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.78"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.404"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "3.899.06"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -4.44%  "
$ws.Range("D14").Value = "66.454.95"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.56%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.366.07"
$ws.Range("E16").Value = "  +1.19%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000164"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "431.16"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  -2.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "3.446.88"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("E26").Value = "  -0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.199"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.60%  "
$ws.Range("E28").Value = "  -3.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.04"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.44"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.60%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("E34").Value = "  -1.68%  "
$ws.Range("E35").Value = "  -2.94%  "
$ws.Range("E36").Value = "  -2.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.16%  "
$ws.Range("E38").Value = "  -2.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").Value = "2.889.33"
$ws.Range("E40").Value = "  +1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.768"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.83%  "
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "40.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0668"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.09%  "
$ws.Range("E46").Value = "  -2.94%  "
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "318.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.90%  "
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.105"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.62%  "
